$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 153.29167
$ws.Range("I33").Value = 155.65218
$ws.Range("K33").Value = 155.65218
$ws.Range("M33").Value = 73.34782000000001
$ws.Range("H58").Value = 2907.6
$ws.Range("I58").Value = 519.5
$ws.Range("K58").Value = 1558.5
$ws.Range("M58").Value = -1408.5
$ws.Range("H86").Value = 4094.6667
$ws.Range("J86").Value = 4226.1665
$ws.Range("L86").Value = 4226.1665
$ws.Range("N86").Value = -6472.1665
$ws.Range("H89").Value = 4094.6667
$ws.Range("J89").Value = 4226.1665
$ws.Range("L89").Value = 21130.8325
$ws.Range("N89").Value = -32362.8325
$ws.Range("H107").Value = 862.4545000000001
$ws.Range("J107").Value = 993.8182
$ws.Range("L107").Value = 993.8182
$ws.Range("N107").Value = -4833.8182
$ws.Range("H131").Value = 5212.769
$ws.Range("J131").Value = 13799.8
$ws.Range("L131").Value = 41399.39999999999
$ws.Range("N131").Value = -51479.39999999999
$ws.Range("H135").Value = 555.6667
$ws.Range("I135").Value = 441.6
$ws.Range("K135").Value = 3974.4
$ws.Range("M135").Value = -1439.4
$ws.Range("H138").Value = 4273.5386
$ws.Range("J138").Value = 4411.1
$ws.Range("L138").Value = 13233.3
$ws.Range("N138").Value = -23513.3
$ws.Range("H141").Value = 1847.8462
$ws.Range("I141").Value = 1823.5
$ws.Range("K141").Value = 5470.5
$ws.Range("M141").Value = -290.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H25").Value = 6003.75
$ws.Range("I25").Value = 4671.6665
$ws.Range("K25").Value = 4671.6665
$ws.Range("M25").Value = -4269.6665
$ws.Range("H32").Value = 20004.164
$ws.Range("I32").Value = 24531.238
$ws.Range("J32").Value = 5378.231
$ws.Range("K32").Value = 24531.238
$ws.Range("L32").Value = 5378.231
$ws.Range("M32").Value = -24244.238
$ws.Range("N32").Value = -5952.231
$ws.Range("H74").Value = 2756.1875
$ws.Range("I74").Value = 2391.1538
$ws.Range("J74").Value = 4338
$ws.Range("K74").Value = 2391.1538
$ws.Range("L74").Value = 4338
$ws.Range("M74").Value = -1517.1538
$ws.Range("N74").Value = -6086
$ws.Range("H77").Value = 2756.1875
$ws.Range("I77").Value = 2391.1538
$ws.Range("J77").Value = 4338
$ws.Range("K77").Value = 11955.769
$ws.Range("L77").Value = 21690
$ws.Range("M77").Value = -7587.769
$ws.Range("N77").Value = -30426
$ws.Range("H132").Value = 41173.117
$ws.Range("I132").Value = 47400.59
$ws.Range("K132").Value = 142201.77
$ws.Range("M132").Value = -139671.77

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 96891.73
$ws.Range("I20").Value = 148259.86
$ws.Range("K20").Value = 148259.86
$ws.Range("M20").Value = -148012.86
$ws.Range("H64").Value = 1712.9
$ws.Range("I64").Value = 1992.8
$ws.Range("J64").Value = 1433
$ws.Range("K64").Value = 1992.8
$ws.Range("L64").Value = 1433
$ws.Range("M64").Value = -1767.8
$ws.Range("N64").Value = -1883
$ws.Range("H67").Value = 1712.9
$ws.Range("I67").Value = 1992.8
$ws.Range("J67").Value = 1433
$ws.Range("K67").Value = 1992.8
$ws.Range("L67").Value = 1433
$ws.Range("M67").Value = -1212.8
$ws.Range("N67").Value = -2993
$ws.Range("H80").Value = 600.5714
$ws.Range("I80").Value = 623
$ws.Range("K80").Value = 623
$ws.Range("M80").Value = 375
$ws.Range("H83").Value = 600.5714
$ws.Range("I83").Value = 623
$ws.Range("K83").Value = 3115
$ws.Range("M83").Value = 1877
$ws.Range("H107").Value = 2245.3635
$ws.Range("I107").Value = 1164.5
$ws.Range("J107").Value = 3542.4
$ws.Range("K107").Value = 1164.5
$ws.Range("L107").Value = 3542.4
$ws.Range("M107").Value = 755.5
$ws.Range("N107").Value = -7382.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5999
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5999
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5999
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -6223
$ws.Range("H31").Value = 2053.2273
$ws.Range("I31").Value = 2094.9473
$ws.Range("K31").Value = 2094.9473
$ws.Range("M31").Value = -1799.9473
$ws.Range("H34").Value = 2053.2273
$ws.Range("I34").Value = 2094.9473
$ws.Range("K34").Value = 2094.9473
$ws.Range("M34").Value = -1892.9473
$ws.Range("H58").Value = 114379.11
$ws.Range("I58").Value = 114379.11
$ws.Range("K58").Value = 114379.11
$ws.Range("M58").Value = -114176.11
$ws.Range("H136").Value = 114379.11
$ws.Range("I136").Value = 114379.11
$ws.Range("K136").Value = 343137.33
$ws.Range("M136").Value = -340587.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 370743.6
$ws.Range("I4").Value = 391317.4
$ws.Range("J4").Value = 415
$ws.Range("K4").Value = 1173952.2
$ws.Range("L4").Value = 1245
$ws.Range("M4").Value = -1173840.2
$ws.Range("N4").Value = -1469
$ws.Range("H7").Value = 2082.2222
$ws.Range("I7").Value = 332.5
$ws.Range("K7").Value = 997.5
$ws.Range("M7").Value = -885.5
$ws.Range("H62").Value = 9249.5
$ws.Range("J62").Value = 9249.5
$ws.Range("L62").Value = 27748.5
$ws.Range("N62").Value = -29120.5
$ws.Range("H65").Value = 9249.5
$ws.Range("J65").Value = 9249.5
$ws.Range("L65").Value = 83245.5
$ws.Range("N65").Value = -90109.5
$ws.Range("H107").Value = 323713.2
$ws.Range("J107").Value = 358311.1
$ws.Range("L107").Value = 1074933.3
$ws.Range("N107").Value = -1078773.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4782
$ws.Range("I70").Value = 4495.2
$ws.Range("K70").Value = 4495.2
$ws.Range("M70").Value = -4225.2
$ws.Range("H73").Value = 4782
$ws.Range("I73").Value = 4495.2
$ws.Range("K73").Value = 4495.2
$ws.Range("M73").Value = -3559.2
$ws.Range("H128").Value = 108499.5
$ws.Range("J128").Value = 108499.5
$ws.Range("L128").Value = 108499.5
$ws.Range("N128").Value = -118459.5
$ws.Range("H132").Value = 58670.5
$ws.Range("I132").Value = 65504.438
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 196513.314
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -193983.314
$ws.Range("N132").Value = -17057
$ws.Range("H141").Value = 55535.285
$ws.Range("J141").Value = 55535.285
$ws.Range("L141").Value = 55535.285
$ws.Range("N141").Value = -65895.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 186831
$ws.Range("I22").Value = 223299.2
$ws.Range("J22").Value = 4490
$ws.Range("K22").Value = 223299.2
$ws.Range("L22").Value = 4490
$ws.Range("M22").Value = -223004.2
$ws.Range("N22").Value = -5080
$ws.Range("H27").Value = 186831
$ws.Range("I27").Value = 223299.2
$ws.Range("J27").Value = 4490
$ws.Range("K27").Value = 223299.2
$ws.Range("L27").Value = 4490
$ws.Range("M27").Value = -223192.2
$ws.Range("N27").Value = -4704
$ws.Range("H46").Value = 24789.6
$ws.Range("J46").Value = 5949.5
$ws.Range("L46").Value = 5949.5
$ws.Range("N46").Value = -6325.5
$ws.Range("H48").Value = 29999.5
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -51321
$ws.Range("H132").Value = 77711.94
$ws.Range("I132").Value = 88035.92999999999
$ws.Range("K132").Value = 264107.79
$ws.Range("M132").Value = -261577.79

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 32512.166
$ws.Range("I43").Value = 10026.333
$ws.Range("K43").Value = 10026.333
$ws.Range("M43").Value = -9877.333000000001
$ws.Range("H95").Value = 30333.334
$ws.Range("J95").Value = 30333.334
$ws.Range("L95").Value = 30333.334
$ws.Range("N95").Value = -35825.334
$ws.Range("H132").Value = 27868.871
$ws.Range("I132").Value = 28544.37
$ws.Range("J132").Value = 2200
$ws.Range("K132").Value = 85633.11
$ws.Range("L132").Value = 6600
$ws.Range("M132").Value = -83103.11
$ws.Range("N132").Value = -11660
